$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# The "Content Placeholder 2" shape currently reads (paragraph numbers):
#   6: "Fundamentals.csv " + "- SEC 10K filings for each company  1.3MB"
#   7: "" (blank line)
#   8: "This data was used in a SQL database..."
#
# Insert a new paragraph ("Total Data size is 103.6MB") right after
# paragraph 6, pushing the existing blank paragraph (and everything after
# it) down by one - leaving the blank line intact between the new
# paragraph and the "This data was used..." paragraph.
$fundamentals = $tr.Paragraphs(6, 1)
$newPara = $fundamentals.InsertAfter("`rTotal Data size is 103.6MB")

# Re-fetch the text range so paragraph indices reflect the inserted text.
$tr = $sh.TextFrame.TextRange
$totalPara = $tr.Paragraphs(7, 1)

$run1 = $totalPara.Characters(1, 6)
$run1.Font.Size = 24

$run2 = $totalPara.Characters(7, $totalPara.Length - 6)
$run2.Font.Size = 24
